# Refresh the cryptos price/volume table (columns D = Price, E = Volume(1h))
# with the latest scraped figures for rows 2-51.
#
# D-column values are strings that look like plain numbers (e.g. "235.50",
# "0.0781", "37.765.50" which uses "." as both a thousands separator and
# decimal point). Assigning such literals straight to .Value lets Excel's
# COM layer auto-coerce them into real numbers, silently dropping the
# original text formatting. To keep them as text we flip the cell to the
# "@" (Text) number format before writing, then reset the style back to
# "Normal" afterwards so no stray style index lingers on a cell that was
# unstyled before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '37.765.50'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.66%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.113.40'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +1.89%  '
$ws.Range("E4").Value = '  -0.02%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '235.50'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("E6").Value = '  +0.35%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '58.28'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").Value = '  -0.03%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.392'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +1.25%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.0781'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.18%  '
$ws.Range("E11").Value = '  +0.94%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '2.422.68'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.73%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '14.63'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("E14").Value = '  +0.54%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.788'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.00%  '
$ws.Range("E16").Value = '  +0.54%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.115.13'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.18%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '37.676.14'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.58%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.21'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.16%  '
$ws.Range("E20").Value = '  +0.78%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.0₃0823'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.80%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '227.62'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("E25").Value = '  -2.48%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '168.10'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.98%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '8.97'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("E28").Value = '  +3.53%  '
$ws.Range("E29").Value = '  -4.34%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '19.48'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.18%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.118'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.33%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.64'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +2.81%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0622'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.74%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.57'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.50%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '4.59'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.33%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '3.49'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +4.33%  '
$ws.Range("E37").Value = '  +0.80%  '
$ws.Range("E38").Value = '  -0.11%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '5.65'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -6.75%  '
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("E41").Value = '  +1.41%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '97.79'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +2.28%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.476.39'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.77%  '
$ws.Range("E44").Value = '  +0.83%  '
$ws.Range("E45").Value = '  -1.00%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '4.22'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -10.56%  '
$ws.Range("E47").Value = '  +1.42%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '15.62'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.91%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '3.04'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +3.42%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '7.33'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.88%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '2.307.14'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.73%  '
